$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Username was not being created because the value stored was "Admin"
# (capitalized) while the expected login username is lowercase "admin".
$ws.Range("A2").Value = "admin"

# Move the active selection to A3 (below the data) instead of D2.
$ws.Range("A3").Select()
